# semana 46 de 2025
# Update weekly case counts ("casos") for several weeks ("semana") in the
# eda2018 sheet, per the updated surveillance data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = 311   # semana 1
$ws.Range("B3").Value  = 755   # semana 2
$ws.Range("B10").Value = 497   # semana 9
$ws.Range("B11").Value = 574   # semana 10
$ws.Range("B12").Value = 473   # semana 11
$ws.Range("B15").Value = 457   # semana 14
$ws.Range("B16").Value = 425   # semana 15
$ws.Range("B18").Value = 457   # semana 17
$ws.Range("B46").Value = 417   # semana 45
$ws.Range("B47").Value = 438   # semana 46
